$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: LP1912 ----------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 05:44:02"
$ws1.Range("A3").Value = "Total filas: 17"

$ws1.Cells.Item(6, 1).Value = "05:44:02"
$ws1.Cells.Item(6, 2).Value = "05:44"
$ws1.Cells.Item(6, 3).Value = "14_ABASTO"
$ws1.Cells.Item(6, 4).Value = 0
$ws1.Cells.Item(6, 5).Value = "LP1912"
$ws1.Cells.Item(7, 1).Value = "05:44:02"
$ws1.Cells.Item(7, 2).Value = "05:47"
$ws1.Cells.Item(7, 3).Value = "17_ROMERO"
$ws1.Cells.Item(7, 4).Value = 3
$ws1.Cells.Item(7, 5).Value = "LP1912"
$ws1.Cells.Item(8, 1).Value = "05:44:02"
$ws1.Cells.Item(8, 2).Value = "06:09"
$ws1.Cells.Item(8, 3).Value = "10_OLMOS"
$ws1.Cells.Item(8, 4).Value = 25
$ws1.Cells.Item(8, 5).Value = "LP1912"
$ws1.Cells.Item(9, 1).Value = "05:44:02"
$ws1.Cells.Item(9, 2).Value = "06:16"
$ws1.Cells.Item(9, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(9, 4).Value = 32
$ws1.Cells.Item(9, 5).Value = "LP1912"
$ws1.Cells.Item(10, 1).Value = "05:44:02"
$ws1.Cells.Item(10, 2).Value = "06:30"
$ws1.Cells.Item(10, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(10, 4).Value = 46
$ws1.Cells.Item(10, 5).Value = "LP1912"
$ws1.Cells.Item(11, 1).Value = "05:44:02"
$ws1.Cells.Item(11, 2).Value = "06:34"
$ws1.Cells.Item(11, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(11, 4).Value = 50
$ws1.Cells.Item(11, 5).Value = "LP1912"
$ws1.Cells.Item(12, 1).Value = "05:44:02"
$ws1.Cells.Item(12, 2).Value = "06:40"
$ws1.Cells.Item(12, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(12, 4).Value = 56
$ws1.Cells.Item(12, 5).Value = "LP1912"
$ws1.Cells.Item(13, 1).Value = "05:44:02"
$ws1.Cells.Item(13, 2).Value = "06:41"
$ws1.Cells.Item(13, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(13, 4).Value = 57
$ws1.Cells.Item(13, 5).Value = "LP1912"
$ws1.Cells.Item(14, 1).Value = "05:44:02"
$ws1.Cells.Item(14, 2).Value = "06:57"
$ws1.Cells.Item(14, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(14, 4).Value = 73
$ws1.Cells.Item(14, 5).Value = "LP1912"
$ws1.Cells.Item(15, 1).Value = "05:44:02"
$ws1.Cells.Item(15, 2).Value = "06:59"
$ws1.Cells.Item(15, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(15, 4).Value = 75
$ws1.Cells.Item(15, 5).Value = "LP1912"
$ws1.Cells.Item(16, 1).Value = "05:44:02"
$ws1.Cells.Item(16, 2).Value = "07:16"
$ws1.Cells.Item(16, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(16, 4).Value = 92
$ws1.Cells.Item(16, 5).Value = "LP1912"
$ws1.Cells.Item(17, 1).Value = "05:44:02"
$ws1.Cells.Item(17, 2).Value = "07:19"
$ws1.Cells.Item(17, 3).Value = "14_ABASTO"
$ws1.Cells.Item(17, 4).Value = 95
$ws1.Cells.Item(17, 5).Value = "LP1912"
$ws1.Cells.Item(18, 1).Value = "05:44:02"
$ws1.Cells.Item(18, 2).Value = "07:21"
$ws1.Cells.Item(18, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(18, 4).Value = 97
$ws1.Cells.Item(18, 5).Value = "LP1912"
$ws1.Cells.Item(19, 1).Value = "05:44:02"
$ws1.Cells.Item(19, 2).Value = "07:22"
$ws1.Cells.Item(19, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(19, 4).Value = 98
$ws1.Cells.Item(19, 5).Value = "LP1912"
$ws1.Cells.Item(20, 1).Value = "05:44:02"
$ws1.Cells.Item(20, 2).Value = "07:29"
$ws1.Cells.Item(20, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(20, 4).Value = 105
$ws1.Cells.Item(20, 5).Value = "LP1912"
$ws1.Cells.Item(21, 1).Value = "05:44:02"
$ws1.Cells.Item(21, 2).Value = "07:35"
$ws1.Cells.Item(21, 3).Value = "10_OLMOS"
$ws1.Cells.Item(21, 4).Value = 111
$ws1.Cells.Item(21, 5).Value = "LP1912"
$ws1.Cells.Item(22, 1).Value = "05:44:02"
$ws1.Cells.Item(22, 2).Value = "07:37"
$ws1.Cells.Item(22, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(22, 4).Value = 113
$ws1.Cells.Item(22, 5).Value = "LP1912"

# ---------- Sheet 2: LP1912-215 ----------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 05:44:02"

$ws2.Cells.Item(6, 1).Value = "05:44:02"
$ws2.Cells.Item(6, 2).Value = "06:16"
$ws2.Cells.Item(6, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(6, 4).Value = 32
$ws2.Cells.Item(6, 5).Value = "LP1912"
$ws2.Cells.Item(7, 1).Value = "05:44:02"
$ws2.Cells.Item(7, 2).Value = "06:57"
$ws2.Cells.Item(7, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(7, 4).Value = 73
$ws2.Cells.Item(7, 5).Value = "LP1912"
$ws2.Cells.Item(8, 1).Value = "05:44:02"
$ws2.Cells.Item(8, 2).Value = "07:16"
$ws2.Cells.Item(8, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(8, 4).Value = 92
$ws2.Cells.Item(8, 5).Value = "LP1912"

# ---------- Sheet 3: 6203-6173 ----------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 05:44:02"
$ws3.Range("A3").Value = "Total filas: 1"

$ws3.Cells.Item(5, 1).Value = "Hora_Scrap"
$ws3.Cells.Item(5, 2).Value = "Hora_Llegada"
$ws3.Cells.Item(5, 3).Value = "Linea"
$ws3.Cells.Item(5, 4).Value = "Minutos"
$ws3.Cells.Item(5, 5).Value = "Parada"

$ws3.Cells.Item(6, 1).Value = "05:44:02"
$ws3.Cells.Item(6, 2).Value = "07:43"
$ws3.Cells.Item(6, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6, 4).Value = 119
$ws3.Cells.Item(6, 5).Value = "L6173"
